$d = $word.ActiveDocument

# 1) "AGULA WALTER DE LA CRUZ" -> "BRUNO JOEL AGUIRRE ENRIQUEZ " (cover page Autor line)
$d.Paragraphs.Item(12).Range.Text = "BRUNO JOEL AGUIRRE ENRIQUEZ "

# 2) Cover page "Asesor:" name -> new advisor name
$d.Paragraphs.Item(14).Range.Text = "AUGUSTO RICARDO MORENO RODRIGUEZ"

# 3) ORCID code: last digit 8 -> 6
$d.Paragraphs.Item(15).Range.Text = "Codigo ORCID: https://orcid.org/0000-0003-3388-4346"

# 4) "2. AUTOR" section author name -> "BRUNO JOEL AGUIRRE ENRIQUEZ "
$d.Paragraphs.Item(25).Range.Text = "BRUNO JOEL AGUIRRE ENRIQUEZ "

# 5) "3. ASESOR" section advisor name -> ". AUGUSTO RICARDO MORENO RODRIGUEZ"
$d.Paragraphs.Item(27).Range.Text = ". AUGUSTO RICARDO MORENO RODRIGUEZ"

# 6) Clear the "NOMBRADO" and "Contabilidad y Finanzas" lines (now empty paragraphs)
$d.Paragraphs.Item(28).Range.Text = ""
$d.Paragraphs.Item(29).Range.Text = ""

# 7) Remove the trailing "none" paragraph entirely
$d.Paragraphs.Item(30).Range.Delete()
